$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.075347494238984
$ws.Range("D2").Value = 1.074379477558459
$ws.Range("E2").Value = 1.078507792967301
$ws.Range("F2").Value = 1.087576891009373
$ws.Range("I2").Value = 1.046326922860092
$ws.Range("J2").Value = 1.08025241187971
$ws.Range("K2").Value = 1.077068845344237
$ws.Range("L2").Value = 1.081186282913324
$ws.Range("M2").Value = 1.090231797435665

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.077035760507895
$ws.Range("D3").Value = 1.075711047986397
$ws.Range("E3").Value = 1.079997171536835
$ws.Range("F3").Value = 1.089111649988355
$ws.Range("I3").Value = 1.046708783593401
$ws.Range("J3").Value = 1.081596655948687
$ws.Range("K3").Value = 1.078216266905026
$ws.Range("L3").Value = 1.082491921752738
$ws.Range("M3").Value = 1.09158443268365

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.078126363441433
$ws.Range("D4").Value = 1.076570790845502
$ws.Range("E4").Value = 1.080959382884948
$ws.Range("F4").Value = 1.090103217238502
$ws.Range("I4").Value = 1.046953668928188
$ws.Range("J4").Value = 1.082464277786314
$ws.Range("K4").Value = 1.078956307803066
$ws.Range("L4").Value = 1.08333472572961
$ws.Range("M4").Value = 1.092457643699516

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.07858442836955
$ws.Range("D5").Value = 1.076931786110254
$ws.Range("E5").Value = 1.081363542114098
$ws.Range("F5").Value = 1.090519715901063
$ws.Range("I5").Value = 1.047056093813832
$ws.Range("J5").Value = 1.082828509117821
$ws.Range("K5").Value = 1.079266849327885
$ws.Range("L5").Value = 1.083688562069232
$ws.Range("M5").Value = 1.092824262657341

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.078661314901922
$ws.Range("D6").Value = 1.076992373184147
$ws.Range("E6").Value = 1.081431381691098
$ws.Range("F6").Value = 1.090589627223857
$ws.Range("I6").Value = 1.047073260701551
$ws.Range("J6").Value = 1.082889635092002
$ws.Range("K6").Value = 1.079318957298889
$ws.Range("L6").Value = 1.083747944920359
$ws.Range("M6").Value = 1.092885791741685

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.07813248578069
$ws.Range("D7").Value = 1.076575616202675
$ws.Range("E7").Value = 1.080964784661053
$ws.Range("F7").Value = 1.090108783901476
$ws.Range("I7").Value = 1.046955039593173
$ws.Range("J7").Value = 1.082469146682579
$ws.Range("K7").Value = 1.078960459510462
$ws.Range("L7").Value = 1.083339455577029
$ws.Range("M7").Value = 1.092462544350638

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.075918434692867
$ws.Range("D8").Value = 1.074829880262403
$ws.Range("E8").Value = 1.079011454528688
$ws.Range("F8").Value = 1.088095890916566
$ws.Range("I8").Value = 1.046456432520006
$ws.Range("J8").Value = 1.080707164773379
$ws.Range("K8").Value = 1.077457126699342
$ws.Range("L8").Value = 1.081627955129655
$ws.Range("M8").Value = 1.090689352983392

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.072002585176691
$ws.Range("D9").Value = 1.071738989729478
$ws.Range("E9").Value = 1.075557440223386
$ws.Range("F9").Value = 1.084536847490895
$ws.Range("I9").Value = 1.045560818173206
$ws.Range("J9").Value = 1.077585148644583
$ws.Range("K9").Value = 1.074789226581033
$ws.Range("L9").Value = 1.078596155686683
$ws.Range("M9").Value = 1.087548817532607

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.069381653465447
$ws.Range("D10").Value = 1.06966804516771
$ws.Range("E10").Value = 1.073246157110561
$ws.Range("F10").Value = 1.082155484694596
$ws.Range("I10").Value = 1.044952133548823
$ws.Range("J10").Value = 1.075491726565732
$ws.Range("K10").Value = 1.07299751269485
$ws.Range("L10").Value = 1.076563752887087
$ws.Range("M10").Value = 1.085443894222536

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.068244151090148
$ws.Range("D11").Value = 1.068768740526491
$ws.Range("E11").Value = 1.072243183490085
$ws.Range("F11").Value = 1.081122148797084
$ws.Range("I11").Value = 1.044685773592075
$ws.Range("J11").Value = 1.074582268767458
$ws.Range("K11").Value = 1.072218469546641
$ws.Range("L11").Value = 1.075680930490795
$ws.Range("M11").Value = 1.08452966018626

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.067821225061397
$ws.Range("D12").Value = 1.068434302959435
$ws.Range("E12").Value = 1.071870297480896
$ws.Range("F12").Value = 1.080737982049439
$ws.Range("I12").Value = 1.044586412189103
$ws.Range("J12").Value = 1.074243995829714
$ws.Range("K12").Value = 1.071928605804063
$ws.Range("L12").Value = 1.075352583735223
$ws.Range("M12").Value = 1.084189644070813

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.067911962759787
$ws.Range("D13").Value = 1.068506059014197
$ws.Range("E13").Value = 1.071950298213142
$ws.Range("F13").Value = 1.080820402670567
$ws.Range("I13").Value = 1.044607744765545
$ws.Range("J13").Value = 1.074316577496773
$ws.Range("K13").Value = 1.071990804986172
$ws.Range("L13").Value = 1.075423034718709
$ws.Range("M13").Value = 1.084262598250924

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.06820920027169
$ws.Range("D14").Value = 1.068741103941221
$ws.Range("E14").Value = 1.072212367551064
$ws.Range("F14").Value = 1.081090400447948
$ws.Range("I14").Value = 1.044677569006027
$ws.Range("J14").Value = 1.074554316460791
$ws.Range("K14").Value = 1.072194519415774
$ws.Range("L14").Value = 1.075653798017288
$ws.Range("M14").Value = 1.084501563155276

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.068392283786813
$ws.Range("D15").Value = 1.068885870151445
$ws.Range("E15").Value = 1.072373792154339
$ws.Range("F15").Value = 1.08125670964143
$ws.Range("I15").Value = 1.04472053380894
$ws.Range("J15").Value = 1.074700734043718
$ws.Range("K15").Value = 1.072319969083384
$ws.Range("L15").Value = 1.075795921998692
$ws.Range("M15").Value = 1.084648740188022

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.069457089398971
$ws.Range("D16").Value = 1.069727674068464
$ws.Range("E16").Value = 1.07331267449337
$ws.Range("F16").Value = 1.08222401673691
$ws.Range("I16").Value = 1.044969751789146
$ws.Range("J16").Value = 1.075552020367406
$ws.Range("K16").Value = 1.073049146650224
$ws.Range("L16").Value = 1.076622283529604
$ws.Range("M16").Value = 1.085504509328318

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.070124302707439
$ws.Range("D17").Value = 1.070255020245219
$ws.Range("E17").Value = 1.073901021472305
$ws.Range("F17").Value = 1.082830189013911
$ws.Range("I17").Value = 1.045125328775188
$ws.Range("J17").Value = 1.076085201463149
$ws.Range("K17").Value = 1.073505672425775
$ws.Range("L17").Value = 1.077139887519965
$ws.Range("M17").Value = 1.086040557169958

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.070513224883639
$ws.Range("D18").Value = 1.070562364827341
$ws.Range("E18").Value = 1.074243985384741
$ws.Range("F18").Value = 1.083183548466898
$ws.Range("I18").Value = 1.045215804709953
$ws.Range("J18").Value = 1.076395908968274
$ws.Range("K18").Value = 1.07377164617727
$ws.Range("L18").Value = 1.077441529742712
$ws.Range("M18").Value = 1.086352956449703

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.070645794936827
$ws.Range("D19").Value = 1.070667119685485
$ws.Range("E19").Value = 1.074360892094791
$ws.Range("F19").Value = 1.083303999514139
$ws.Range("I19").Value = 1.04524660905974
$ws.Range("J19").Value = 1.076501803656193
$ws.Range("K19").Value = 1.073862284034353
$ws.Range("L19").Value = 1.077544336878932
$ws.Range("M19").Value = 1.086459431312998

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.070052743200398
$ws.Range("D20").Value = 1.070198466661256
$ws.Range("E20").Value = 1.073837919095829
$ws.Range("F20").Value = 1.082765174324465
$ws.Range("I20").Value = 1.045108664736536
$ws.Range("J20").Value = 1.076028026041801
$ws.Range("K20").Value = 1.07345672365985
$ws.Range("L20").Value = 1.077084381242275
$ws.Range("M20").Value = 1.085983072139223

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.068121682574053
$ws.Range("D21").Value = 1.068671900069764
$ws.Range("E21").Value = 1.07213520404931
$ws.Range("F21").Value = 1.081010902298604
$ws.Range("I21").Value = 1.044657019224952
$ws.Range("J21").Value = 1.074484321037545
$ws.Range("K21").Value = 1.072134544258118
$ws.Range("L21").Value = 1.0755858558494
$ws.Range("M21").Value = 1.084431205866922

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.066905186755777
$ws.Range("D22").Value = 1.067709792487276
$ws.Range("E22").Value = 1.071062684757617
$ws.Range("F22").Value = 1.079905950555697
$ws.Range("I22").Value = 1.044370600043725
$ws.Range("J22").Value = 1.073511066912109
$ws.Range("K22").Value = 1.071300384131476
$ws.Range("L22").Value = 1.074641196853811
$ws.Range("M22").Value = 1.083452999528733

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.067550302447521
$ws.Range("D23").Value = 1.068220044570811
$ws.Range("E23").Value = 1.071631436350124
$ws.Range("F23").Value = 1.080491896850685
$ws.Range("I23").Value = 1.044522669830758
$ws.Range("J23").Value = 1.074027263451748
$ws.Range("K23").Value = 1.071742861624719
$ws.Range("L23").Value = 1.075142216495939
$ws.Range("M23").Value = 1.083971804252531

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.07008507864028
$ws.Range("D24").Value = 1.070224021556484
$ws.Range("E24").Value = 1.073866432989238
$ws.Range("F24").Value = 1.08279455231491
$ws.Range("I24").Value = 1.045116195331064
$ws.Range("J24").Value = 1.07605386204363
$ws.Range("K24").Value = 1.073478842458202
$ws.Range("L24").Value = 1.077109462966353
$ws.Range("M24").Value = 1.08600904797951

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.073016703076372
$ws.Range("D25").Value = 1.07253984795603
$ws.Range("E25").Value = 1.076451862814411
$ws.Range("F25").Value = 1.085458432361665
$ws.Range("I25").Value = 1.04579438906082
$ws.Range("J25").Value = 1.078394354501357
$ws.Range("K25").Value = 1.075481220890121
$ws.Range("L25").Value = 1.079381886122286
$ws.Range("M25").Value = 1.088362663312961

